$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate rows 2-32 (language/count entries) into rows 33-63
for ($i = 2; $i -le 32; $i++) {
    $destRow = $i + 31
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($i, 1).Value2
}
